$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "2021-10-29"
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 17000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 17500
$ws.Range("N2").Value = "`$/saco 25 kilos"
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 700
$ws.Range("Q2").Value = 25

# Row 3
$ws.Range("D3").Value = "2020-11-26"
$ws.Range("J3").Value = 33
$ws.Range("K3").Value = 19000
$ws.Range("L3").Value = 19500
$ws.Range("M3").Value = 19303
$ws.Range("O3").Value = "Provincia de Diguillín"
$ws.Range("P3").Value = 772

# Row 4
$ws.Range("D4").Value = "2020-11-26"
$ws.Range("J4").Value = 34
$ws.Range("K4").Value = 19500
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = 19735
$ws.Range("P4").Value = 789

# Row 5
$ws.Range("D5").Value = "2020-11-27"
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 18500
$ws.Range("L5").Value = 19000
$ws.Range("M5").Value = 18820
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 753

# Row 6
$ws.Range("D6").Value = "2021-11-04"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 17500
$ws.Range("P6").Value = 700

# Row 7
$ws.Range("D7").Value = "2021-03-05"
$ws.Range("J7").Value = 33
$ws.Range("M7").Value = 22545
$ws.Range("P7").Value = 22545

# Row 8
$ws.Range("D8").Value = "2021-11-16"
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15500
$ws.Range("P8").Value = 620

# Row 9
$ws.Range("D9").Value = "2021-01-27"
$ws.Range("K9").Value = 26000
$ws.Range("L9").Value = 28000
$ws.Range("M9").Value = 27048
$ws.Range("N9").Value = "`$/saco 30 kilos"
$ws.Range("P9").Value = 27048
$ws.Range("Q9").Value = 1

# Row 10
$ws.Range("D10").Value = "2021-11-24"
$ws.Range("J10").Value = 80
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 13000
$ws.Range("M10").Value = 12500
$ws.Range("N10").Value = "`$/saco 25 kilos"
$ws.Range("O10").Value = "Región del Maule"
$ws.Range("P10").Value = 500
$ws.Range("Q10").Value = 25

# Row 11
$ws.Range("D11").Value = "2021-03-03"
$ws.Range("J11").Value = 32
$ws.Range("K11").Value = 22000
$ws.Range("L11").Value = 23000
$ws.Range("M11").Value = 22562
$ws.Range("N11").Value = "`$/saco 30 kilos"
$ws.Range("O11").Value = "Región de La Araucanía"
$ws.Range("P11").Value = 22562
$ws.Range("Q11").Value = 1

# Row 12
$ws.Range("D12").Value = "2020-12-09"
$ws.Range("H12").Value = "Perfection"
$ws.Range("J12").Value = 30
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = 19500
$ws.Range("O12").Value = "Región de Ñuble"
$ws.Range("P12").Value = 780

# Row 13
$ws.Range("D13").Value = "2021-11-11"
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 15500
$ws.Range("P13").Value = 620

# Row 15
$ws.Range("D15").Value = "2021-01-20"
$ws.Range("H15").Value = "Perfection"
$ws.Range("J15").Value = 43
$ws.Range("K15").Value = 24000
$ws.Range("L15").Value = 25000
$ws.Range("M15").Value = 24419
$ws.Range("N15").Value = "`$/saco 30 kilos"
$ws.Range("O15").Value = "Región de La Araucanía"
$ws.Range("P15").Value = 24419
$ws.Range("Q15").Value = 1

# Row 16
$ws.Range("D16").Value = "2020-12-11"
$ws.Range("H16").Value = "Perfection"
$ws.Range("J16").Value = 30
$ws.Range("L16").Value = 21000
$ws.Range("M16").Value = 20500
$ws.Range("O16").Value = "Provincia de Diguillín"
$ws.Range("P16").Value = 820

# Row 17
$ws.Range("D17").Value = "2020-11-30"
$ws.Range("J17").Value = 42
$ws.Range("K17").Value = 18000
$ws.Range("L17").Value = 19000
$ws.Range("M17").Value = 18595
$ws.Range("P17").Value = 744

# Row 18
$ws.Range("D18").Value = "2021-11-30"
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 14000
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = 14500
$ws.Range("O18").Value = "Provincia de Diguillín"
$ws.Range("P18").Value = 580

# Row 19
$ws.Range("D19").Value = "2021-11-30"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("K19").Value = 12000
$ws.Range("L19").Value = 13000
$ws.Range("M19").Value = 12500
$ws.Range("P19").Value = 500

# Row 20
$ws.Range("D20").Value = "2020-11-24"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("J20").Value = 28
$ws.Range("L20").Value = 19500
$ws.Range("M20").Value = 19268
$ws.Range("O20").Value = "Provincia de Diguillín"
$ws.Range("P20").Value = 771

# Row 21
$ws.Range("D21").Value = "2020-11-24"
$ws.Range("J21").Value = 56
$ws.Range("K21").Value = 19000
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = 19464
$ws.Range("P21").Value = 779

# Row 22
$ws.Range("D22").Value = "2021-11-29"
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = 12000
$ws.Range("L22").Value = 13000
$ws.Range("M22").Value = 12500
$ws.Range("N22").Value = "`$/saco 25 kilos"
$ws.Range("O22").Value = "Región del Maule"
$ws.Range("P22").Value = 500
$ws.Range("Q22").Value = 25

# Row 23
$ws.Range("D23").Value = "2021-11-18"
$ws.Range("K23").Value = 15000
$ws.Range("L23").Value = 16000
$ws.Range("M23").Value = 15500
$ws.Range("P23").Value = 620

# Row 24
$ws.Range("D24").Value = "2021-11-19"
$ws.Range("J24").Value = 60

# Row 25
$ws.Range("D25").Value = "2021-01-19"
$ws.Range("J25").Value = 42
$ws.Range("K25").Value = 19000
$ws.Range("L25").Value = 20000
$ws.Range("M25").Value = 19524
$ws.Range("O25").Value = "Región de La Araucanía"
$ws.Range("P25").Value = 781

# Row 26
$ws.Range("D26").Value = "2020-12-01"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("J26").Value = 48
$ws.Range("K26").Value = 17000
$ws.Range("L26").Value = 18000
$ws.Range("M26").Value = 17479
$ws.Range("N26").Value = "`$/saco 25 kilos"
$ws.Range("O26").Value = "Región del Maule"
$ws.Range("P26").Value = 699
$ws.Range("Q26").Value = 25

# Row 27
$ws.Range("D27").Value = "2021-11-22"
$ws.Range("J27").Value = 60
$ws.Range("K27").Value = 13000
$ws.Range("L27").Value = 14000
$ws.Range("M27").Value = 13500
$ws.Range("P27").Value = 540

# Row 28
$ws.Range("D28").Value = "2021-03-12"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("J28").Value = 43
$ws.Range("K28").Value = 24000
$ws.Range("L28").Value = 25000
$ws.Range("M28").Value = 24581
$ws.Range("N28").Value = "`$/saco 30 kilos"
$ws.Range("O28").Value = "Región de La Araucanía"
$ws.Range("P28").Value = 24581
$ws.Range("Q28").Value = 1

# Row 29
$ws.Range("D29").Value = "2021-11-25"
$ws.Range("J29").Value = 80
$ws.Range("K29").Value = 12000
$ws.Range("L29").Value = 13000
$ws.Range("M29").Value = 12500
$ws.Range("P29").Value = 500

# Row 30
$ws.Range("D30").Value = "2020-11-25"
$ws.Range("J30").Value = 62
$ws.Range("K30").Value = 19000
$ws.Range("L30").Value = 20000
$ws.Range("M30").Value = 19516
$ws.Range("P30").Value = 781

# Row 31
$ws.Range("D31").Value = "2021-03-18"
$ws.Range("J31").Value = 22
$ws.Range("K31").Value = 20000
$ws.Range("L31").Value = 22000
$ws.Range("M31").Value = 21091
$ws.Range("O31").Value = "Región de La Araucanía"
$ws.Range("P31").Value = 844

# Row 32
$ws.Range("D32").Value = "2020-12-02"
$ws.Range("H32").Value = "Perfection"
$ws.Range("J32").Value = 40
$ws.Range("K32").Value = 18000
$ws.Range("L32").Value = 19000
$ws.Range("M32").Value = 18500
$ws.Range("O32").Value = "Región del Maule"
$ws.Range("P32").Value = 740

# Row 33
$ws.Range("D33").Value = "2021-11-09"
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 17000
$ws.Range("L33").Value = 18000
$ws.Range("M33").Value = 17500
$ws.Range("P33").Value = 700
